$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.459.95'
$ws.Range('E2').Value = '  -0.79%  '
$ws.Range('D3').Value = '1.634.70'
$ws.Range('E3').Value = '  -1.25%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  +0.38%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.000'
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '303.60'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3786'
$ws.Range('E7').Value = '  +0.45%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '51.96'
$ws.Range('E8').Value = '  -0.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3625'
$ws.Range('E9').Value = '  -0.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08182'
$ws.Range('E10').Value = '  +0.29%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.229'
$ws.Range('E11').Value = '  -3.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9995'
$ws.Range('E12').Value = '  +0.37%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.40'
$ws.Range('E13').Value = '  -3.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.480'
$ws.Range('E14').Value = '  -3.18%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.379'
$ws.Range('E15').Value = '  -0.08%  '
$ws.Range('D17').Value = '1.628.73'
$ws.Range('E17').Value = '  -1.33%  '
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06948'
$ws.Range('E19').Value = '  +0.94%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.570'
$ws.Range('E20').Value = '  -0.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.53'
$ws.Range('E21').Value = '  -4.82%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  +0.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.56'
$ws.Range('D24').Value = '23.461.21'
$ws.Range('E24').Value = '  -0.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.519'
$ws.Range('E25').Value = '  +4.22%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.064'
$ws.Range('E26').Value = '  -3.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.14'
$ws.Range('E27').Value = '  -1.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '151.13'
$ws.Range('E28').Value = '  +0.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.284'
$ws.Range('E29').Value = '  -0.69%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.06'
$ws.Range('E30').Value = '  -2.73%  '
$ws.Range('D31').Value = '1.811.90'
$ws.Range('E31').Value = '  -1.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.641'
$ws.Range('E32').Value = '  -3.61%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.157'
$ws.Range('E33').Value = '  -6.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.048'
$ws.Range('E34').Value = '  +7.73%  '
$ws.Range('E35').Value = '  +2.92%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02760'
$ws.Range('E36').Value = '  -2.56%  '
$ws.Range('B37').Value = 'Algorand'
$ws.Range('C37').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2496'
$ws.Range('E37').Value = '  -2.61%  '
$ws.Range('B38').Value = 'Stellar'
$ws.Range('C38').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.08782'
$ws.Range('E38').Value = '  -1.25%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.07111'
$ws.Range('E39').Value = '  -3.75%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.023'
$ws.Range('E40').Value = '  -5.45%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.7005'
$ws.Range('E41').Value = '  -1.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.342'
$ws.Range('E42').Value = '  -2.72%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '15.94'
$ws.Range('E43').Value = '  -2.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '12.15'
$ws.Range('E44').Value = '  -3.60%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6496'
$ws.Range('E45').Value = '  -1.66%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.0000'
$ws.Range('E46').Value = '  +0.45%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.274'
$ws.Range('E47').Value = '  -3.67%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.970'
$ws.Range('E48').Value = '  -1.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07972'
$ws.Range('E49').Value = '  -0.84%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '126.77'
$ws.Range('E50').Value = '  -2.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.188'
$ws.Range('E51').Value = '  -2.66%  '
